$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.57%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.65"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.243"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.44%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05694"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.45%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.12%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.187"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.24%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8498"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.77%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8541"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.93%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1369"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.40%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07087"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.09%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03153"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "7.74%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09203"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.96%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001538"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.78%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005947"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.87%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006027"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.51%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.491"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.15%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.174"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.35%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3167"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.52%"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03257"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.62%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.09%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.487"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.59%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04079"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.20%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1378"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.11%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001220"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.08%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004140"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.52%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.87%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03753"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.23%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.80%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.75%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.50%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009326"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.74%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005264"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.98%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07496"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "15.87%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.01%"
